$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are formatted as text so numeric-looking strings are preserved exactly
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","B8","C8","D8","E8","B9","C9","D9","E9","B10","C10","D10","E10","B11","C11","D11","E11","B12","C12","D12","E12","B13","C13","D13","E13","B14","C14","D14","E14","B15","C15","D15","E15","B16","C16","D16","E16","B17","C17","D17","E17","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","D48","E48")
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "317.27"
$ws.Range("E2").Value = "3.61%"
$ws.Range("D3").Value = "39.77"
$ws.Range("E3").Value = "1.12%"
$ws.Range("D4").Value = "5.153"
$ws.Range("E4").Value = "0.75%"
$ws.Range("D5").Value = "0.08244"
$ws.Range("E5").Value = "2.13%"
$ws.Range("D6").Value = "2.048"
$ws.Range("E6").Value = "6.24%"
$ws.Range("D7").Value = "8.353"
$ws.Range("E7").Value = "3.79%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.322"
$ws.Range("E8").Value = "2.92%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9396"
$ws.Range("E9").Value = "1.48%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1353"
$ws.Range("E10").Value = "-3.13%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1995"
$ws.Range("E11").Value = "4.08%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09087"
$ws.Range("E12").Value = "0.31%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03502"
$ws.Range("E13").Value = "-0.60%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09805"
$ws.Range("E14").Value = "0.27%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001408"
$ws.Range("E15").Value = "0.99%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006145"
$ws.Range("E16").Value = "4.91%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.687"
$ws.Range("E17").Value = "-2.05%"
$ws.Range("E18").Value = "-1.40%"
$ws.Range("D19").Value = "0.3494"
$ws.Range("E19").Value = "0.94%"
$ws.Range("D20").Value = "0.1322"
$ws.Range("E20").Value = "-0.24%"
$ws.Range("D21").Value = "4.962"
$ws.Range("E21").Value = "6.03%"
$ws.Range("D22").Value = "0.2452"
$ws.Range("E22").Value = "1.52%"
$ws.Range("D23").Value = "0.04356"
$ws.Range("E23").Value = "-0.34%"
$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").Value = "2.65%"
$ws.Range("D25").Value = "0.004797"
$ws.Range("E25").Value = "12.06%"
$ws.Range("E26").Value = "-0.03%"
$ws.Range("D39").Value = "0.02303"
$ws.Range("E39").Value = "13.13%"
$ws.Range("D40").Value = "0.05186"
$ws.Range("E40").Value = "3.11%"
$ws.Range("D41").Value = "0.007755"
$ws.Range("E41").Value = "3.02%"
$ws.Range("E42").Value = "7.96%"
$ws.Range("D43").Value = "0.1407"
$ws.Range("E43").Value = "4.80%"
$ws.Range("D44").Value = "0.002083"
$ws.Range("E44").Value = "-0.49%"
$ws.Range("D45").Value = "0.009313"
$ws.Range("E45").Value = "-4.95%"
$ws.Range("D46").Value = "0.00006923"
$ws.Range("E46").Value = "11.57%"
$ws.Range("D48").Value = "0.002887"
$ws.Range("E48").Value = "0.15%"
